$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Date Added" column (column C)
# so the new layout becomes: Pack, Price, Cost, Date Added, Is Sold
$ws.Columns.Item(3).Insert()

# Header for the new column
$ws.Range("C1").Value = "Cost"

# Cost values for each data row
$ws.Range("C2").Value = 17
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 7.77
$ws.Range("C5").Value = 46.86
$ws.Range("C6").Value = 3
